$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = ' A Megawatt (MW) measures power, or the rate of energy consumption per unit time. On the other hand, a Megawatt-hour (MWh) is a measure of energy, specifically, one megawatt being sustained for one hour. In simple terms, 1 MW is the amount of power needed to light 1000 100-watt light bulbs continuously for an hour. So, 1 MWh represents the energy equivalent of running 1 MW for one hour.'
$ws.Range("C3").Value = ' The Annual Energy Outlook was not specifically mentioned in the provided text. It is a report published by the U.S. Energy Information Administration (EIA), which provides projections of energy supply, demand, and prices for the next 25-30 years.'
$ws.Range("C4").Value = ' The Clean Air Act is a legislation in the United States that requires the U.S. Environmental Protection Agency to set National Ambient Air Quality Standards for six common air pollutants. These pollutants are particle pollution, ground-level ozone, carbon monoxide, sulfur oxides, nitrogen oxides, and lead.'
$ws.Range("C5").Value = ' An Electric Generating Unit (EGU) refers to any facility that produces electricity through the conversion of energy from various sources, such as fossil fuels, nuclear power, or renewable resources like solar, wind, ocean thermal, wave action, and tidal action. In the context provided, it specifically refers to a facility burning fossil fuel and is used in discussions related to reducing emissions and improving air quality.'
$ws.Range("C6").Value = ' Electricity Dispatch Models (EDMs) optimize the dispatch of a power system by determining the optimal operating conditions for generating units in real-time or near real-time to meet demand and ensure grid stability. They consider multiple factors, including:
1. Load forecasts: Accurate predictions of electricity consumption are essential for proper dispatch. Load forecasting accounts for weather conditions, time of day, holidays, and other factors that can affect energy demand.
2. Generation availability: EDMs take into account the current operational status of generators, including maintenance schedules, fuel supply, and outages.
3. Generator characteristics: The unique technical specifications of each generating unit, such as minimum/maximum output limits, ramp rates, and efficiency levels, are also considered by EDMs to ensure that they operate within their optimal range.
4. Market pricing signals: Real-time energy market prices can affect the economic dispatch of power plants. For example, generators may prioritize running when prices are high or curtail production when prices are low.
5. Grid constraints: EDMs must also consider transmission line limitations, such as congestion, voltage limits, and capacity factors to ensure that power flows efficiently and safely within the grid.
6. Emission targets: To minimize emissions, EDMs can prioritize dispatching cleaner generating units, like renewable energy sources, and optimize their output.
7. Operating reserve requirements: Reserve margins are necessary to maintain grid reliability and stability in case of unexpected events such as generator outages or sudden changes in load. EDMs must ensure there is an adequate level of reserves available for contingencies.
In summary, EDMs optimize the dispatch of a power system by balancing electricity demand with available generation while considering a range of operational, market, and grid constraints to maintain reliability, stability, and emissions compliance.'
$ws.Range("C7").Value = ' The Emissions & Generation Resource Integrated Database is maintained by the U.S. Environmental Protection Agency (EPA).'
$ws.Range("C8").Value = ' Examples of Energy Efficiency policies, programs and projects include:
1. Building codes that require energy-efficient construction materials and designs for new buildings.
2. Appliance standards that set minimum efficiency requirements for appliances such as refrigerators, air conditioners, and lighting fixtures.
3. Lighting retrofit programs that replace inefficient incandescent bulbs with more efficient compact fluorescent or LED bulbs.
4. Energy-efficient product labeling programs that help consumers make informed choices about energy-efficient products.
5. Utility demand-response programs that reduce electricity use during peak demand periods by offering financial incentives to customers who reduce their electricity consumption.
6. Industrial process optimization programs that identify and implement energy efficiency measures in manufacturing processes.
7. Energy-efficient transportation programs such as electric vehicles, carpooling, and bike sharing initiatives.
8. Smart grid technologies that optimize the use of electricity across the power grid to improve efficiency and reduce waste.'
$ws.Range("C9").Value = ' In the context provided, there is no definition or mention of "ISOs" (Interconnection System Operators). However, if you''re referring to International Organization for Standardization (ISO), it''s an independent, non-governmental international organization with a membership of 164 national standards bodies. ISO is responsible for developing and publishing a wide range of standards related to various industries and sectors. These standards aim to provide guidelines and specifications for products, services, and systems, ensuring quality, safety, efficiency, and compatibility.'
$ws.Range("C10").Value = ' Based on the provided text, there is no direct mention of "Power Pools" in the context given. Power pools are usually regional organizations that coordinate and manage electricity generation and transmission across multiple utilities to maintain a stable power grid. They operate on the principle of sharing resources and balancing supply and demand. However, without additional context or information, I cannot confidently say if this term is relevant to the text provided.'
$ws.Range("C11").Value = ' The Public Utilities Commission (PUC) or Public Service Commission (PSC), depending on the state, is a regulatory body that oversees utilities such as electricity, natural gas, water, and sometimes telecommunications. They are responsible for setting rates, rules, and regulations to ensure fair prices, reliable service, safety, and sustainability within the energy sector. In some cases, they may also play a role in promoting renewable energy through policies like Renewable Portfolio Standards or implementing energy efficiency programs. However, specific roles and responsibilities can vary by state.'
$ws.Range("C12").Value = ' Electricity Dispatch Models simulate the operation of a power system, including the generation, transmission, and distribution of electricity from various energy sources such as wind, solar, fossil fuels, etc., to meet the demand of consumers at any given time. These models take into account factors like weather conditions (for renewable energy sources), fuel prices, maintenance schedules for power plants, and demand forecasts. They help in optimizing the use of available resources, minimizing operational costs, maintaining grid stability, and meeting regulatory requirements such as Renewable Portfolio Standards.'
